$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '29.184.94'
$ws.Cells.Item(2, 5).Value = '  +0.31%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.888.98'
$ws.Cells.Item(3, 5).Value = '  -0.44%  '

# Row 4
$ws.Cells.Item(4, 4).Value = "'1.002"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.04%  '

# Row 5
$ws.Cells.Item(5, 4).Value = "'322.24"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -2.95%  '

# Row 6
$ws.Cells.Item(6, 4).Value = "'1.001"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.01%  '

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.4701"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +2.33%  '

# Row 8
$ws.Cells.Item(8, 4).Value = "'0.4023"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -2.59%  '

# Row 9
$ws.Cells.Item(9, 4).Value = "'47.26"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -1.14%  '

# Row 10
$ws.Cells.Item(10, 4).Value = "'0.07999"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +0.07%  '

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.9935"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -1.75%  '

# Row 12
$ws.Cells.Item(12, 4).Value = "'22.67"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +2.29%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '1.886.54'
$ws.Cells.Item(13, 5).Value = '  -0.74%  '

# Row 14
$ws.Cells.Item(14, 4).Value = "'5.905"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -0.54%  '

# Row 15
$ws.Cells.Item(15, 4).Value = "'7.017"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -1.52%  '

# Row 16
$ws.Cells.Item(16, 4).Value = "'89.79"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +0.86%  '

# Row 17
$ws.Cells.Item(17, 4).Value = "'1.002"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.10%  '

# Row 18
$ws.Cells.Item(18, 4).Value = "'0.06631"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +1.02%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -0.74%  '

# Row 20
$ws.Cells.Item(20, 4).Value = "'17.43"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -1.05%  '

# Row 21
$ws.Cells.Item(21, 4).Value = "'1.000"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -0.13%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '29.186.27'
$ws.Cells.Item(22, 5).Value = '  +0.49%  '

# Row 23
$ws.Cells.Item(23, 4).Value = "'5.486"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -0.17%  '

# Row 24
$ws.Cells.Item(24, 4).Value = "'11.64"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +1.94%  '

# Row 25
$ws.Cells.Item(25, 4).Value = "'2.176"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -1.03%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '2.087.53'
$ws.Cells.Item(26, 5).Value = '  -1.63%  '

# Row 27
$ws.Cells.Item(27, 4).Value = "'155.19"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -1.11%  '

# Row 28
$ws.Cells.Item(28, 4).Value = "'19.64"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -0.50%  '

# Row 29
$ws.Cells.Item(29, 4).Value = "'6.022"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +6.88%  '

# Row 30
$ws.Cells.Item(30, 4).Value = "'2.079"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -1.96%  '

# Row 31
$ws.Cells.Item(31, 4).Value = "'118.62"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +1.31%  '

# Row 32
$ws.Cells.Item(32, 4).Value = "'1.027"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -2.01%  '

# Row 33
$ws.Cells.Item(33, 4).Value = "'0.09408"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +0.30%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'HuobiToken'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(34, 4).Value = "'3.540"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +0.17%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'ARBITRUM'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(35, 4).Value = "'1.380"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -2.59%  '

# Row 36
$ws.Cells.Item(36, 4).Value = "'5.355"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +0.06%  '

# Row 37
$ws.Cells.Item(37, 4).Value = "'0.06052"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -0.68%  '

# Row 38
$ws.Cells.Item(38, 4).Value = "'0.02234"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -0.19%  '

# Row 39
$ws.Cells.Item(39, 4).Value = "'1.172"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -0.47%  '

# Row 40
$ws.Cells.Item(40, 4).Value = "'7.994"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -5.21%  '

# Row 41
$ws.Cells.Item(41, 4).Value = "'0.5813"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -0.32%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Algorand'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(42, 4).Value = "'0.1827"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -0.08%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'RenderToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(43, 4).Value = "'2.451"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +6.46%  '

# Row 44
$ws.Cells.Item(44, 4).Value = "'10.03"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -0.99%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'WEMIXToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(45, 4).Value = "'1.273"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +2.15%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'Cronos'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(46, 4).Value = "'0.07697"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +2.59%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  +0.20%  '

# Row 48
$ws.Cells.Item(48, 4).Value = "'0.5465"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -1.09%  '

# Row 49
$ws.Cells.Item(49, 4).Value = "'1.900"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -1.18%  '

# Row 50
$ws.Cells.Item(50, 4).Value = "'113.13"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.83%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'WOONetwork'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Cells.Item(51, 4).Value = "'0.2940"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +0.35%  '
